$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("week4")

# --- week4 ("Week 4") sheet: fill in the new progress-report entries ---
# (order matches the original authoring order so new shared-string indices
# line up: 74..81)
$ws4.Range("C1").Value  = "Cuddy Esmeres"
$ws4.Range("C2").Value  = "Xymon Malquisto"

$ws4.Range("A7").Value  = "Finish code logic for polar bar chart "
$ws4.Range("C7").Value  = "Done"

$ws4.Range("A8").Value  = "calculator"

$ws4.Range("A10").Value = "Finish ~80% of logic code"
$ws4.Range("C10").Value = "Done"

$ws4.Range("A14").Value = "Coding, debugging, researching"

$ws4.Range("A28").Value = "Properly track time usage to allow "
$ws4.Range("A29").Value = "focusing on other tasks needed to be"
$ws4.Range("A30").Value = "finished"

$ws4.Range("A21").Value = "Effiiency in code creating"

# --- switch the active/selected tab from week3 to week4 ---
$ws4.Activate()
$ws4.Range("C2").Select()
